# Auto-generated edit script: updates crypto price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-CellText "D2" "26.300.35"
Set-CellText "E2" "  +0.02%  "

Set-CellText "D3" "1.601.97"
Set-CellText "E3" "  +0.45%  "

Set-CellText "D4" "0.999"
Set-CellText "E4" "  -0.01%  "

Set-CellText "D5" "213.22"
Set-CellText "E5" "  +0.19%  "

Set-CellText "D6" "0.502"
Set-CellText "E6" "  -0.11%  "

Set-CellText "E7" "  +0.01%  "

Set-CellText "E8" "  -0.54%  "

Set-CellText "D9" "0.0608"
Set-CellText "E9" "  -0.20%  "

Set-CellText "D10" "19.08"
Set-CellText "E10" "  -1.57%  "

Set-CellText "E11" "  +0.30%  "

Set-CellText "D12" "1.823.16"
Set-CellText "E12" "  +0.34%  "

Set-CellText "D13" "1.588.19"
Set-CellText "E13" "  -0.47%  "

Set-CellText "D14" "4.03"
Set-CellText "E14" "  -0.22%  "

Set-CellText "D15" "0.510"
Set-CellText "E15" "  -2.26%  "

Set-CellText "D16" "64.07"
Set-CellText "E16" "  -0.60%  "

Set-CellText "D17" "26.305.67"
Set-CellText "E17" "  +0.09%  "

Set-CellText "D18" "0.0₃0724"
Set-CellText "E18" "  -0.49%  "

Set-CellText "D19" "7.48"

Set-CellText "D20" "216.75"
Set-CellText "E20" "  +1.42%  "

Set-CellText "E21" "  +0.00%  "

Set-CellText "D22" "4.33"
Set-CellText "E22" "  +0.90%  "

Set-CellText "D23" "9.02"
Set-CellText "E23" "  +0.00%  "

Set-CellText "E24" "  -1.64%  "

Set-CellText "D25" "144.77"
Set-CellText "E25" "  +0.09%  "

Set-CellText "E26" "  +0.04%  "

Set-CellText "D27" "7.00"
Set-CellText "E27" "  -0.89%  "

Set-CellText "E28" "  +0.90%  "

Set-CellText "D29" "15.21"
Set-CellText "E29" "  -0.14%  "

Set-CellText "E30" "  -0.41%  "

Set-CellText "E31" "  +0.21%  "

Set-CellText "E32" "  +0.09%  "

Set-CellText "D33" "1.434.84"
Set-CellText "E33" "  +7.11%  "

Set-CellText "D34" "2.98"
Set-CellText "E34" "  +0.79%  "

Set-CellText "D35" "2.41"
Set-CellText "E35" "  -1.17%  "

Set-CellText "E36" "  -0.26%  "

Set-CellText "D37" "0.556"
Set-CellText "E37" "  -5.82%  "

Set-CellText "D38" "0.0166"
Set-CellText "E38" "  -0.31%  "

Set-CellText "D39" "0.830"
Set-CellText "E39" "  +0.97%  "

Set-CellText "D40" "5.78"
Set-CellText "E40" "  +1.28%  "

Set-CellText "E41" "  +0.06%  "

Set-CellText "E42" "  +1.60%  "

Set-CellText "B43" "RocketPoolETH"
Set-CellText "C43" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText "D43" "1.737.73"
Set-CellText "E43" "  +0.60%  "

Set-CellText "B44" "TrustWalletToken"
Set-CellText "C44" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D44" "0.761"
Set-CellText "E44" "  -0.43%  "

Set-CellText "E45" "  -11.67%  "

Set-CellText "D46" "61.14"
Set-CellText "E46" "  -1.37%  "

Set-CellText "D47" "87.19"
Set-CellText "E47" "  +2.17%  "

Set-CellText "B48" "BabyDogeCoin"
Set-CellText "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText "D48" "0.0₆0103"
Set-CellText "E48" "  -1.46%  "

Set-CellText "B49" "RenderToken"
Set-CellText "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D49" "1.48"
Set-CellText "E49" "  -0.52%  "

Set-CellText "B50" "Cronos"
Set-CellText "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText "D50" "0.0501"
Set-CellText "E50" "  -0.29%  "

Set-CellText "B51" "Algorand"
Set-CellText "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-CellText "D51" "0.0957"
Set-CellText "E51" "  -1.92%  "
